$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh: two brand-new daily price rows are inserted right
# after the header (pushing every existing record down by two rows), and
# the table grows from A1:R54 to A1:R56.
$ws.Rows("2:3").Insert()

# Inherit the regular data-row formatting (incl. the date number format on
# column D) from the row immediately below instead of the header styling
# that Insert() copies by default.
$ws.Range("A4:R4").Copy()
$ws.Range("A2:R3").PasteSpecial(-4122)

# Row 2: new "Calameño" / "Segunda" price entry
$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(2,3).Value = "Arica y Parinacota"
$ws.Cells.Item(2,4).Value = 44599
$ws.Cells.Item(2,5).Value = 15
$ws.Cells.Item(2,6).Value = 100112027
$ws.Cells.Item(2,7).Value = "Melón"
$ws.Cells.Item(2,8).Value = "Calameño"
$ws.Cells.Item(2,9).Value = "Segunda"
$ws.Cells.Item(2,10).Value = 30
$ws.Cells.Item(2,11).Value = 6000
$ws.Cells.Item(2,12).Value = 6000
$ws.Cells.Item(2,13).Value = 6000
$ws.Cells.Item(2,14).Value = "$/caja 24 unidades"
$ws.Cells.Item(2,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(2,16).Value = 250
$ws.Cells.Item(2,17).Value = 24
$ws.Cells.Item(2,18).Value = "Hortaliza"

# Row 3: new "Tuna" / "Segunda" price entry
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(3,3).Value = "Arica y Parinacota"
$ws.Cells.Item(3,4).Value = 44599
$ws.Cells.Item(3,5).Value = 15
$ws.Cells.Item(3,6).Value = 100112027
$ws.Cells.Item(3,7).Value = "Melón"
$ws.Cells.Item(3,8).Value = "Tuna"
$ws.Cells.Item(3,9).Value = "Segunda"
$ws.Cells.Item(3,10).Value = 30
$ws.Cells.Item(3,11).Value = 5000
$ws.Cells.Item(3,12).Value = 5000
$ws.Cells.Item(3,13).Value = 5000
$ws.Cells.Item(3,14).Value = "$/caja 24 unidades"
$ws.Cells.Item(3,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(3,16).Value = 208
$ws.Cells.Item(3,17).Value = 24
$ws.Cells.Item(3,18).Value = "Hortaliza"
